$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.232.80'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.604.55'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''212.58'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("D10").Value = '''18.44'
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '1.827.98'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.603.92'
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '26.186.47'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '''61.96'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''200.84'
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("D25").Value = '''143.94'
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("D28").Value = '''15.20'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("E30").Value = '  +3.86%  '
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  +2.41%  '
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("D35").Value = '''2.38'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '1.162.08'
$ws.Range("E36").Value = '  +4.24%  '
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -0.38%  '
$ws.Range("D40").Value = '''0.785'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.33'
$ws.Range("E42").Value = '  +4.13%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.783'
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").Value = '1.739.65'
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").Value = '0.0₆0105'
$ws.Range("E46").Value = '  +18.24%  '
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").Value = '''54.07'
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("E51").Value = '  -0.09%  '
